# Fix average stars calculation - now per attack not per war
#
# 1. ROSTER: Krunal (row 22) and LittleSinn (row 23) had their
#    "Wars Participated" counts corrected now that attacks are
#    counted per-attack instead of per-war.
# 2. War sheet "20260119T102025-000Z": Krunal's first attack (row 31),
#    which was previously recorded as a miss (all zeros), actually
#    happened - 3 stars / 100% destruction.
# 3. Same sheet: LittleSinn made a second attack that was missing from
#    the report entirely. A new row is inserted for it (after her first
#    attack in row 39), pushing every following row down by one.

$wb = $excel.ActiveWorkbook

# --- 1. ROSTER sheet -------------------------------------------------
$roster = $wb.Worksheets.Item("ROSTER")
$roster.Range("E22").Value = 1
$roster.Range("E23").Value = 4

# --- 2 & 3. War sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("20260119T102025-000Z")

# Krunal's first (only) attack actually landed - correct the stats.
$ws.Range("J31").Value = 1
$ws.Range("K31").Value = 3
$ws.Range("L31").Value = 100
$ws.Range("M31").Value = "Yes"
$ws.Range("N31").Value = "No"

# Insert a new row for LittleSinn's second attack, shifting rows
# 40-61 down to 41-62.
$ws.Rows("40:40").Insert()

$ws.Range("A40").Value = "20260119T102025-000Z"
$ws.Range("B40").Value = "inWar"
$ws.Range("C40").Value = "No"
$ws.Range("D40").Value = "20260121T092025.000Z"
$ws.Range("E40").Value = 35
$ws.Range("F40").Value = "LittleSinn"
$ws.Range("G40").Value = "#YGV99UU"
$ws.Range("H40").Value = 17
$ws.Range("I40").Value = 30
$ws.Range("J40").Value = 2
$ws.Range("K40").Value = 2
$ws.Range("L40").Value = 91
$ws.Range("M40").Value = "No"
$ws.Range("N40").Value = "No"
$ws.Range("O40").Value = "No"
